# Add harvester ("S.GISH") to every data row, an experiment design
# ("90minuteInduction") to every data row, and record the strain used for
# the first sample ("KN99alpha").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

$ws.Cells.Item(2, 6).Value = "KN99alpha"

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Match the author's on-screen selection when the workbook was saved.
$ws.Range("B2:B11").Select()
